$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 292
    $ws.Range("F3").Value = 234
    $ws.Range("F4").Value = 37
}
